$wb = $excel.ActiveWorkbook

# The plate layout uses columns C..W (spans 8 cols -> repeats every 12) to hold the
# "elution control" well labels for rows 4-15. They used to be packed into
# consecutive rows (4,5,6,7) which doesn't line up with 8-channel pipetting.
# Re-space them onto every second row (4,6,8,10) instead.
$sheetNames = @("merged", "sample_name", "experimental_unit", "condition", "replicate")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Copy (not Value2 read+write) so the shared-string cell type survives
    # intact even when the text looks numeric (e.g. the literal "1"), then
    # Clear the vacated source cell outright so it drops out of the sheet
    # entirely instead of lingering as a styled-but-empty cell. Move
    # back-to-front so each source is fully drained before it is reused as a
    # destination.
    $ws.Range("C7").Copy($ws.Range("C10"))
    $ws.Range("C7").Clear()

    $ws.Range("C6").Copy($ws.Range("C8"))
    $ws.Range("C6").Clear()

    $ws.Range("C5").Copy($ws.Range("C6"))
    $ws.Range("C5").Clear()
}

# Update the recorded selections / active sheet to match the new layout.
$ws1 = $wb.Worksheets.Item("merged")
$ws1.Activate()
$ws1.Range("AF7").Select()

$ws2 = $wb.Worksheets.Item("sample_name")
$ws2.Activate()
$ws2.Range("C12").Select()

$ws3 = $wb.Worksheets.Item("experimental_unit")
$ws3.Activate()
$ws3.Range("C10").Select()

$ws4 = $wb.Worksheets.Item("condition")
$ws4.Activate()
$ws4.Range("C5").Select()

$ws5 = $wb.Worksheets.Item("replicate")
$ws5.Activate()
$ws5.Range("AF6").Select()
